$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: IS_TINATETT_PRODUCT / Yes
$ws.Range("F1").Value = "IS_TINATETT_PRODUCT"
$ws.Range("F2").Value = "Yes"

# New column widths (for F and G)
$ws.Columns.Item(6).ColumnWidth = 19.25
$ws.Columns.Item(7).ColumnWidth = 23.42

# Update view: zoom and selection
$ws.Application.ActiveWindow.Zoom = 125
$ws.Range("F13").Select()
